$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Links de noticias (investing.com) para cada ativo da base, coletados via scraping.
$newsLinks = @{
    3 = "https://www.investing.com/equities/microsoft-corp-news"
    4 = "https://www.investing.com/equities/google-inc-news"
    5 = "https://www.investing.com/equities/amazon-com-inc-news"
    6 = "https://www.investing.com/equities/nvidia-corp-news"
    7 = "https://www.investing.com/equities/meta-platforms-inc-news"
    8 = "https://www.investing.com/equities/tesla-motors-news"
    9 = "https://www.investing.com/equities/berkshire-hathaway-news"
    10 = "https://www.investing.com/equities/eli-lilly-and-company-news"
    11 = "https://www.investing.com/equities/visa-inc-news"
    12 = "https://www.investing.com/equities/jp-morgan-chase-news"
    13 = "https://www.investing.com/equities/exxon-mobil-corp-news"
    14 = "https://www.investing.com/equities/johnson-johnson-news"
    15 = "https://www.investing.com/equities/mastercard-inc-news"
    16 = "https://www.investing.com/equities/procter-gamble-news"
    17 = "https://www.investing.com/equities/costco-whsl-corp-new-news"
    18 = "https://www.investing.com/equities/bank-of-america-corp-news"
    19 = "https://www.investing.com/equities/netflix-inc-news"
    20 = "https://www.investing.com/equities/adv-micro-device-news"
    21 = "https://www.investing.com/equities/coca-cola-co-news"
    22 = "https://www.investing.com/equities/pepsico-inc-news"
    23 = "https://www.investing.com/equities/wal-mart-stores-news"
    24 = "https://www.investing.com/equities/mcdonalds-corp-news"
    25 = "https://www.investing.com/equities/disney-news"
    26 = "https://www.investing.com/equities/caterpillar-inc-news"
    27 = "https://www.investing.com/equities/intel-corp-news"
    28 = "https://www.investing.com/equities/cisco-sys-inc-news"
    29 = "https://www.investing.com/equities/oracle-corp-news"
    30 = "https://www.investing.com/equities/salesforce-com-news"
    31 = "https://www.investing.com/equities/adobe-sys-inc-news"
    32 = "https://www.investing.com/equities/nike-news"
    33 = "https://www.investing.com/equities/starbucks-corp-news"
    34 = "https://www.investing.com/equities/boeing-co-news"
    35 = "https://www.investing.com/equities/goldman-sachs-group-news"
    36 = "https://www.investing.com/equities/morgan-stanley-news"
    37 = "https://www.investing.com/equities/ford-motor-co-news"
    38 = "https://www.investing.com/equities/gen-motors-news"
    39 = "https://www.investing.com/equities/pfizer-news"
    40 = "https://www.investing.com/equities/chevron-corp-news"
    41 = "https://www.investing.com/equities/paypal-holdings-inc-news"
    42 = "https://www.investing.com/equities/coinbase-global-inc-news"
    43 = "https://www.investing.com/equities/uber-technologies-inc-news"
    44 = "https://www.investing.com/equities/airbnb-inc-news"
}

# Linhas que devem virar hyperlinks "de verdade" (igual ao padrao ja usado na F2 do Apple)
$hyperlinkRows = @(3, 4, 44)

foreach ($row in ($newsLinks.Keys | Sort-Object)) {
    $url = $newsLinks[$row]
    $cell = $ws.Cells.Item($row, 6)
    if ($hyperlinkRows -contains $row) {
        $ws.Hyperlinks.Add($cell, $url)
    } else {
        $cell.Value = $url
    }
}

[void]$ws.Range("F9").Select()
